$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 had redundant "value" labels copy-pasted across C1:F1 - trim back to A1:B1
$ws.Range("C1:F1").ClearContents()

# "Model" parameter row is renamed to "production_function"
$ws.Range("A8").Value = "production_function"

# Insert a new "L_curve" row right after production_function
$ws.Range("A9").EntireRow.Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").NumberFormat = "0.00E+00"
$ws.Range("B9").Value = 1

# The old "Strain"/"Deletion" extra row (now shifted down to row 17) is removed
$ws.Range("A17").EntireRow.Delete()

# This sheet becomes the active / selected sheet, with the last row selected
$ws.Rows(17).Select()
